$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-21 Friday" "2025-11-22 Saturday"

Replace-Text "75÷9=8, 3" "25÷8=3, 1"
Replace-Text "45÷7=6, 3" "85÷8=10, 5"
Replace-Text "67÷4=16, 3" "74÷9=8, 2"
Replace-Text "65÷2=32, 1" "40÷6=6, 4"
Replace-Text "13÷7=1, 6" "89÷5=17, 4"

Replace-Text "33÷2=16, 1" "22÷6=3, 4"
Replace-Text "33÷4=8, 1" "45÷9=5, 0"
Replace-Text "42÷4=10, 2" "23÷2=11, 1"
Replace-Text "41÷2=20, 1" "70÷9=7, 7"
Replace-Text "61÷6=10, 1" "34÷2=17, 0"

Replace-Text "75÷8=9, 3" "22÷9=2, 4"
Replace-Text "25÷5=5, 0" "11÷6=1, 5"
Replace-Text "29÷6=4, 5" "68÷8=8, 4"
Replace-Text "23÷5=4, 3" "21÷8=2, 5"
Replace-Text "84÷3=28, 0" "75÷6=12, 3"

Replace-Text "10÷2=5, 0" "17÷6=2, 5"
Replace-Text "37÷4=9, 1" "18÷3=6, 0"
Replace-Text "56÷4=14, 0" "48÷3=16, 0"
Replace-Text "97÷5=19, 2" "59÷7=8, 3"
Replace-Text "51÷4=12, 3" "19÷2=9, 1"

Replace-Text "83÷8=10, 3" "36÷5=7, 1"
Replace-Text "74÷4=18, 2" "81÷6=13, 3"
Replace-Text "27÷6=4, 3" "76÷2=38, 0"
Replace-Text "86÷6=14, 2" "60÷5=12, 0"
Replace-Text "70÷5=14, 0" "96÷3=32, 0"
